$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new display value. Values that Excel would otherwise
# auto-convert to a number (plain decimals with no thousand-dot
# grouping) are prefixed with a leading apostrophe so they are
# stored as text, matching the source data feed (e.g. "187.66").
$updates = @{
    'D2' = '75.746.58'
    'E2' = '  +8.49%  '
    'D3' = '2.723.06'
    'E3' = '  +11.68%  '
    'E4' = '  +0.11%  '
    'D5' = '''187.66'
    'E5' = '  +12.00%  '
    'D6' = '''591.99'
    'E6' = '  +4.40%  '
    'E7' = '  -0.06%  '
    'E8' = '  +5.33%  '
    'E9' = '  +14.10%  '
    'D10' = '2.726.63'
    'E10' = '  +11.91%  '
    'E11' = '  +1.29%  '
    'E12' = '  +8.34%  '
    'D13' = '''4.80'
    'E13' = '  +1.76%  '
    'D14' = '3.227.00'
    'E14' = '  +11.89%  '
    'D15' = '75.546.55'
    'E15' = '  +8.45%  '
    'D16' = '''0.0000191'
    'E16' = '  +6.52%  '
    'D17' = '''27.13'
    'E17' = '  +12.46%  '
    'D18' = '2.716.42'
    'E18' = '  +11.72%  '
    'D19' = '''9.43'
    'E19' = '  +30.89%  '
    'D20' = '''12.20'
    'E20' = '  +11.92%  '
    'D21' = '''379.24'
    'E21' = '  +10.03%  '
    'D22' = '''2.32'
    'E22' = '  +15.34%  '
    'E23' = '  +6.22%  '
    'E24' = '  +4.48%  '
    'D25' = '''71.35'
    'E25' = '  +7.77%  '
    'D26' = '''0.999'
    'E26' = '  -0.10%  '
    'E27' = '  +10.61%  '
    'D28' = '''9.66'
    'E28' = '  +13.39%  '
    'D29' = '2.858.83'
    'E29' = '  +11.70%  '
    'D30' = '''1.00'
    'E30' = '  -0.27%  '
    'E31' = '  +16.28%  '
    'D32' = '''525.84'
    'E32' = '  +14.95%  '
    'E33' = '  +12.83%  '
    'E34' = '  +6.77%  '
    'E35' = '  +11.16%  '
    'D36' = '''1.00'
    'E36' = '  +0.14%  '
    'D37' = '''0.120'
    'E37' = '  +7.38%  '
    'D38' = '''162.00'
    'E38' = '  +1.70%  '
    'D39' = '''19.64'
    'E39' = '  +7.40%  '
    'D40' = '''19.39'
    'E40' = '  +1.37%  '
    'D42' = '''175.41'
    'E42' = '  +28.62%  '
    'D43' = '''5.08'
    'E43' = '  +14.56%  '
    'D44' = '''1.72'
    'E44' = '  +12.88%  '
    'E45' = '  +9.83%  '
    'D46' = '''1.23'
    'E46' = '  +13.29%  '
    'E47' = '  +14.62%  '
    'D48' = '''39.20'
    'E48' = '  +2.99%  '
    'E49' = '  +18.63%  '
    'E50' = '  +9.62%  '
    'D51' = '''0.551'
    'E51' = '  +11.82%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$wb.Save()
